$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 39 (Karnauf Ábel block): fill in a new task row ---
# (Done before row 6 so new shared-string entries are appended in the same
# order as the source edit: "Bakancs oldal" / "Bakancs oldal készítése"
# before "Minden" / "itt volt minden teso".)
$ws.Range("C3").Copy() | Out-Null
$ws.Range("C39").PasteSpecial(-4122) | Out-Null
$ws.Range("D3").Copy() | Out-Null
$ws.Range("D39").PasteSpecial(-4122) | Out-Null

$ws.Range("B39").Value = "Bakancs oldal"
$ws.Range("C39").Value = (Get-Date -Year 2023 -Month 6 -Day 7 -Hour 18 -Minute 0 -Second 0)
$ws.Range("D39").Value = (Get-Date -Year 2023 -Month 6 -Day 7 -Hour 19 -Minute 0 -Second 0)
$ws.Range("E39").Value = "Kész"
$ws.Range("F39").Value = "Bakancs oldal készítése"

# --- Row 6 (Papp Zsombor block): fill in a new task row ---
# Reuse the existing date-format style (matches other rows in this block)
# by copying the format from a sibling cell before writing the value, so
# Excel doesn't mint a brand-new number format for the date.
$ws.Range("C3").Copy() | Out-Null
$ws.Range("C6").PasteSpecial(-4122) | Out-Null
$ws.Range("D3").Copy() | Out-Null
$ws.Range("D6").PasteSpecial(-4122) | Out-Null

$ws.Range("B6").Value = "Minden"
$ws.Range("C6").Value = (Get-Date -Year 2023 -Month 6 -Day 6 -Hour 22 -Minute 0 -Second 0)
$ws.Range("D6").Value = (Get-Date -Year 2023 -Month 6 -Day 6 -Hour 22 -Minute 30 -Second 0)
$ws.Range("E6").Value = "Kész"
$ws.Range("F6").Value = "itt volt minden teso"

# --- Rows 15-25 (Papp Zsombor block, C/D columns): normalize cell format ---
# D15 and the blank C/D cells below it pick up the same "empty" style used
# elsewhere in the sheet.
$ws.Range("C29").Copy() | Out-Null
$ws.Range("D15").PasteSpecial(-4122) | Out-Null
$ws.Range("C16:D24").PasteSpecial(-4122) | Out-Null
$ws.Range("C37").Copy() | Out-Null
$ws.Range("C25:D25").PasteSpecial(-4122) | Out-Null

# --- Row 28 (Lépesfalvi Balázs block): clear the task row ---
$ws.Range("B28:F28").ClearContents()
$ws.Range("C29").Copy() | Out-Null
$ws.Range("C28").PasteSpecial(-4122) | Out-Null
$ws.Range("D29").Copy() | Out-Null
$ws.Range("D28").PasteSpecial(-4122) | Out-Null

# --- Final UI selection, matching the end state of the edit ---
$ws.Range("G7").Select() | Out-Null
